$wb = $excel.ActiveWorkbook

# Add the new worksheet "Hojita 2" after the existing "Datos" sheet
$datos = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $datos)
$newSheet.Name = "Hojita 2"

# Populate cell A1 with the requested text
$newSheet.Range("A1").Value = "Holaa hojita 2"

# Leave the selection where the author left it, and keep this sheet active
$newSheet.Range("I19").Select()
$newSheet.Activate()
